$d = $word.ActiveDocument
$sec = $d.Sections.First

# The document has two Pearson logo pictures (one in each footer) whose
# inline-shape name was "image2.png" and needs to become "image1.png",
# and one BTEC logo picture (in the "first page" header) whose name was
# "image1.jpg" and needs to become "image2.jpg".

$ftrs = $sec.Footers
for ($i = 1; $i -le $ftrs.Count; $i++) {
    $f = $ftrs.Item($i)
    if ($f.Exists) {
        $cnt = $f.Range.InlineShapes.Count
        for ($j = 1; $j -le $cnt; $j++) {
            $s = $f.Range.InlineShapes.Item($j)
            if ($s.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                $s.Name = "image1.png"
            }
        }
    }
}

$hdrs = $sec.Headers
for ($i = 1; $i -le $hdrs.Count; $i++) {
    $h = $hdrs.Item($i)
    if ($h.Exists) {
        $cnt = $h.Range.InlineShapes.Count
        for ($j = 1; $j -le $cnt; $j++) {
            $s = $h.Range.InlineShapes.Item($j)
            if ($s.AlternativeText -eq "BTec_Logo-Orange") {
                $s.Name = "image2.jpg"
            }
        }
    }
}
